$d = $word.ActiveDocument

# Locate the paragraph that contains "ywkffyencszqptvk" (the last piece of
# real content before the run of trailing blank paragraphs). The diff
# shows the *second* blank paragraph following it (of four) gets turned
# into four new paragraphs of secret-looking content, while the first
# and the last two blank paragraphs are left untouched.
$anchorRng = $d.Content
$anchorRng.Find.Execute("ywkffyencszqptvk", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorStart = $anchorRng.Paragraphs.Item(1).Range.Start

$allParas = $d.Paragraphs
$anchorIndex = -1
for ($i = 1; $i -le $allParas.Count; $i++) {
    if ($allParas.Item($i).Range.Start -eq $anchorStart) {
        $anchorIndex = $i
        break
    }
}

# anchorIndex      -> paragraph with "ywkffyencszqptvk"
# anchorIndex + 1  -> first blank paragraph (stays blank)
# anchorIndex + 2  -> second blank paragraph -> becomes the 4 new paragraphs
$targetIndex = $anchorIndex + 2

$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range

# New paragraph 1: "Vandor Validation token 22.06.25"
$r.Text = "Vandor Validation token 22.06.25"
$r.InsertParagraphAfter()

# New paragraph 2: "Name-" + "Vendor Validator Push" (two runs)
$p2 = $d.Paragraphs.Item($targetIndex + 1)
$r2 = $p2.Range
$r2.Text = "Name-"
$ins2 = $d.Range($r2.End - 1, $r2.End - 1)
$ins2.InsertAfter("Vendor Validator Push")
$p2.Range.InsertParagraphAfter()

# New paragraph 3: "User-" + "KoenigSalary"
$p3 = $d.Paragraphs.Item($targetIndex + 2)
$r3 = $p3.Range
$r3.Text = "User-"
$ins3 = $d.Range($r3.End - 1, $r3.End - 1)
$ins3.InsertAfter("KoenigSalary")
$p3.Range.InsertParagraphAfter()

# New paragraph 4: "Pass and Token -" + "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"
$p4 = $d.Paragraphs.Item($targetIndex + 3)
$r4 = $p4.Range
$r4.Text = "Pass and Token -"
$ins4 = $d.Range($r4.End - 1, $r4.End - 1)
$ins4.InsertAfter("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")

Write-Output "Inserted 4 new paragraphs starting at index $targetIndex"
